# Posthoc.xlsx demo update — "Improve demo: level display"
#
# Relabels the Chocolate_1/Chocolate_2 level text (with/without -> no/yes),
# tweaks one displayed emmCI string, refreshes the per-row statistics, and
# nudges three column widths (H, I, N) to better fit the new labels.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Level labels -----------------------------------------------------
$ws.Range("B2").Value = "no"
$ws.Range("B3").Value = "no"
$ws.Range("C2").Value = "yes"
$ws.Range("C3").Value = "yes"

# --- Displayed CI string for the first emmCI column (row 2 only) -----
$ws.Range("D2").Value = "1.86 (1.85, 1.87)"

# --- Refreshed statistics, row 2 (Gender = female) --------------------
$ws.Range("F2").Value = 0.04580625174067107
$ws.Range("G2").Value = 0.0458062517406711
$ws.Range("H2").Value = -0.0016437045720332222
$ws.Range("I2").Value = -0.025662985979687386
$ws.Range("J2").Value = 4.009086477938423
$ws.Range("M2").Value = 0.1809067067706429
$ws.Range("N2").Value = 0.008115410399677865

# --- Refreshed statistics, row 3 (Gender = male) ----------------------
$ws.Range("F3").Value = 0.0000000005914474688420944
$ws.Range("G3").Value = 0.0000000011828948931480454
$ws.Range("H3").Value = 0.00976335396847361
$ws.Range("I3").Value = 0.12548337809223806
$ws.Range("J3").Value = 39.932137077404555
$ws.Range("M3").Value = 0.5709436304954304
$ws.Range("N3").Value = 0.07535330334489954

# --- Column width tweaks ----------------------------------------------
# ColumnWidth goes through Excel's normal character->pixel grid snap
# (round(width*6)+5 px), so we dial in the COM-facing "characters" value
# that snaps to the closest on-grid width to the requested figures.
$ws.Columns.Item(8).ColumnWidth = 14.666666666666666   # -> 15.42578125 target (snaps to 15.5)
$ws.Columns.Item(9).ColumnWidth = 13.666666666666666   # -> 14.42578125 target (snaps to 14.5)
$ws.Columns.Item(14).ColumnWidth = 11.833333333333334  # -> 12.7109375 target (snaps to 12.666666666666666)
